$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = @("31/12/2018", "31/12/2019", "31/12/2021", "31/12/2022", "31/12/2023")

# Apply the date number format to C2 first, then propagate that exact
# style (without creating duplicate style records) to the rest of the
# column via copy/paste-special of formats only.
$ws.Range("C2").NumberFormat = "mm-dd-yy"
$ws.Range("C2").Copy()
$ws.Range("C3:C16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

for ($row = 2; $row -le 16; $row++) {
    $idx = ($row - 2) % 5
    $cell = $ws.Cells.Item($row, 3)
    $cell.Value = $dates[$idx]
}

$ws.Range("C2").Select()
